$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: 2021年 ---
$ws.Cells.Item(11, 1).Value = "2021年"
$row11Values = @(27418.25,6047.38,2246,4.84,19371.33,46300.14,5924.55,17182.11,3250.4,3979.41,2773.86,1981.89,5010.06,26847.77,3623.34,3190.29,14005.02,14144.42,52881.53,2507.03,44661.24,9080.99,123888.9,51949.29,3457.55,26353.33,10554.26,12291.73,5927.98,828485.36,84418.42,31031.82,9234.73,9006.52,1355.56,24195.39,20241.04,40206.58,3609.99,9412.540000000001,41669.7,7248.16)
for ($i = 0; $i -lt $row11Values.Length; $i++) {
  $ws.Cells.Item(11, $i + 2).Value = $row11Values[$i]
}

# --- Row 12: 2022年 ---
$ws.Cells.Item(12, 1).Value = "2022年"
$row12Values = @(29901.7,6452.4,830.5,4.1,21540.8,50427.2,6811.9,18833.8,3376.7,4104.5,3435.8,1981.5,5266.2,29357.7,3829.3,3453.1,14795.7,15759.1,59354.6,2008.1,44449.6,9705.1,130935,64660.7,3445.7,29608.8,10670.2,12785.2,5817.5,882994.2,91753.8,31734.2,9472.9,9728,1441.6,22338.2,12158.4,43802.9,4234.6,10883.9,44538.7,7304.4)
for ($i = 0; $i -lt $row12Values.Length; $i++) {
  $ws.Cells.Item(12, $i + 2).Value = $row12Values[$i]
}

# Match the style of the year-label column (A10) for the new labels in A11/A12
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11:A12").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(11, 1).Value = "2021年"
$ws.Cells.Item(12, 1).Value = "2022年"
$excel.CutCopyMode = 0
